$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H33").Value = 389.125
$ws_ALC.Range("I33").Value = 187.16667
$ws_ALC.Range("J33").Value = 995
$ws_ALC.Range("K33").Value = 187.16667
$ws_ALC.Range("L33").Value = 995
$ws_ALC.Range("M33").Value = 41.83332999999999
$ws_ALC.Range("N33").Value = -1453

$ws_ALC.Range("H129").Value = 1323.08
$ws_ALC.Range("I129").Value = 551
$ws_ALC.Range("J129").Value = 1837.8
$ws_ALC.Range("K129").Value = 1653
$ws_ALC.Range("L129").Value = 5513.4
$ws_ALC.Range("M129").Value = 3347
$ws_ALC.Range("N129").Value = -15513.4

$ws_ALC.Range("H132").Value = 2912.4773
$ws_ALC.Range("I132").Value = 1866.9688
$ws_ALC.Range("J132").Value = 5700.5
$ws_ALC.Range("K132").Value = 5600.9064
$ws_ALC.Range("L132").Value = 17101.5
$ws_ALC.Range("M132").Value = -3070.9064
$ws_ALC.Range("N132").Value = -22161.5

$ws_ALC.Range("H138").Value = 3371.2964
$ws_ALC.Range("I138").Value = 4500
$ws_ALC.Range("J138").Value = 3114.7727
$ws_ALC.Range("K138").Value = 13500
$ws_ALC.Range("L138").Value = 9344.3181
$ws_ALC.Range("M138").Value = -8360
$ws_ALC.Range("N138").Value = -19624.3181

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 2946.2307
$ws_ARM.Range("I2").Value = 1273.3334
$ws_ARM.Range("J2").Value = 5227.4546
$ws_ARM.Range("K2").Value = 1273.3334
$ws_ARM.Range("L2").Value = 5227.4546
$ws_ARM.Range("M2").Value = -1160.3334
$ws_ARM.Range("N2").Value = -5453.4546

$ws_ARM.Range("H44").Value = 37924.5
$ws_ARM.Range("J44").Value = 37924.5
$ws_ARM.Range("L44").Value = 37924.5
$ws_ARM.Range("N44").Value = -38900.5

$ws_ARM.Range("H80").Value = 27442
$ws_ARM.Range("J80").Value = 27442
$ws_ARM.Range("L80").Value = 27442
$ws_ARM.Range("N80").Value = -29438

$ws_ARM.Range("H83").Value = 27442
$ws_ARM.Range("J83").Value = 27442
$ws_ARM.Range("L83").Value = 82326
$ws_ARM.Range("N83").Value = -92310

$ws_ARM.Range("H88").Value = 2338.375
$ws_ARM.Range("I88").Value = 0
$ws_ARM.Range("J88").Value = 2338.375
$ws_ARM.Range("K88").Value = 0
$ws_ARM.Range("L88").Value = 2338.375
$ws_ARM.Range("N88").Value = -3150.375
$ws_ARM.Range("M88").ClearContents()

$ws_ARM.Range("H91").Value = 2338.375
$ws_ARM.Range("I91").Value = 0
$ws_ARM.Range("J91").Value = 2338.375
$ws_ARM.Range("K91").Value = 0
$ws_ARM.Range("L91").Value = 2338.375
$ws_ARM.Range("N91").Value = -5146.375
$ws_ARM.Range("M91").ClearContents()

$ws_ARM.Range("H116").Value = 2946.2307
$ws_ARM.Range("I116").Value = 1273.3334
$ws_ARM.Range("J116").Value = 5227.4546
$ws_ARM.Range("K116").Value = 1273.3334
$ws_ARM.Range("L116").Value = 5227.4546
$ws_ARM.Range("M116").Value = 1020.6666
$ws_ARM.Range("N116").Value = -9815.454600000001

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 2946.2307
$ws_BSM.Range("I3").Value = 1273.3334
$ws_BSM.Range("J3").Value = 5227.4546
$ws_BSM.Range("K3").Value = 1273.3334
$ws_BSM.Range("L3").Value = 5227.4546
$ws_BSM.Range("M3").Value = -1159.3334
$ws_BSM.Range("N3").Value = -5455.4546

$ws_BSM.Range("H26").Value = 8126.6665
$ws_BSM.Range("I26").Value = 8126.6665
$ws_BSM.Range("J26").Value = 0
$ws_BSM.Range("K26").Value = 8126.6665
$ws_BSM.Range("L26").Value = 0
$ws_BSM.Range("M26").Value = -7834.6665
$ws_BSM.Range("N26").ClearContents()

$ws_BSM.Range("H58").Value = 80000
$ws_BSM.Range("J58").Value = 80000
$ws_BSM.Range("L58").Value = 80000
$ws_BSM.Range("N58").Value = -80588

$ws_BSM.Range("H59").Value = 64925
$ws_BSM.Range("J59").Value = 64925
$ws_BSM.Range("L59").Value = 64925
$ws_BSM.Range("N59").Value = -66619

$ws_BSM.Range("H82").Value = 68837.45
$ws_BSM.Range("I82").Value = 155378.5
$ws_BSM.Range("J82").Value = 19385.428
$ws_BSM.Range("K82").Value = 155378.5
$ws_BSM.Range("L82").Value = 19385.428
$ws_BSM.Range("M82").Value = -154995.5
$ws_BSM.Range("N82").Value = -20151.428

$ws_BSM.Range("H85").Value = 68837.45
$ws_BSM.Range("I85").Value = 155378.5
$ws_BSM.Range("J85").Value = 19385.428
$ws_BSM.Range("K85").Value = 155378.5
$ws_BSM.Range("L85").Value = 19385.428
$ws_BSM.Range("M85").Value = -154052.5
$ws_BSM.Range("N85").Value = -22037.428

$ws_BSM.Range("H126").Value = 46893.332
$ws_BSM.Range("J126").Value = 46893.332
$ws_BSM.Range("L126").Value = 46893.332
$ws_BSM.Range("N126").Value = -56773.332

$ws_BSM.Range("H134").Value = 85493.39
$ws_BSM.Range("I134").Value = 113280.15
$ws_BSM.Range("J134").Value = 2133.111
$ws_BSM.Range("K134").Value = 339840.45
$ws_BSM.Range("L134").Value = 6399.333
$ws_BSM.Range("M134").Value = -337305.45
$ws_BSM.Range("N134").Value = -11469.333

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H20").Value = 47787.09
$ws_CRP.Range("J20").Value = 47787.09
$ws_CRP.Range("L20").Value = 47787.09
$ws_CRP.Range("N20").Value = -48259.09

$ws_CRP.Range("H30").Value = 47787.09
$ws_CRP.Range("J30").Value = 47787.09
$ws_CRP.Range("L30").Value = 47787.09
$ws_CRP.Range("N30").Value = -47969.09

$ws_CRP.Range("H31").Value = 1767.25
$ws_CRP.Range("I31").Value = 1294.1111
$ws_CRP.Range("J31").Value = 3186.6667
$ws_CRP.Range("K31").Value = 1294.1111
$ws_CRP.Range("L31").Value = 3186.6667
$ws_CRP.Range("M31").Value = -999.1111000000001
$ws_CRP.Range("N31").Value = -3776.6667

$ws_CRP.Range("H34").Value = 1767.25
$ws_CRP.Range("I34").Value = 1294.1111
$ws_CRP.Range("J34").Value = 3186.6667
$ws_CRP.Range("K34").Value = 1294.1111
$ws_CRP.Range("L34").Value = 3186.6667
$ws_CRP.Range("M34").Value = -1092.1111
$ws_CRP.Range("N34").Value = -3590.6667

$ws_CRP.Range("H50").Value = 8523.833000000001
$ws_CRP.Range("J50").Value = 8935.091
$ws_CRP.Range("L50").Value = 8935.091
$ws_CRP.Range("N50").Value = -10185.091

$ws_CRP.Range("H109").Value = 20634
$ws_CRP.Range("J109").Value = 20634
$ws_CRP.Range("L109").Value = 20634
$ws_CRP.Range("N109").Value = -22714

$ws_CRP.Range("H128").Value = 47787.09
$ws_CRP.Range("J128").Value = 47787.09
$ws_CRP.Range("L128").Value = 47787.09
$ws_CRP.Range("N128").Value = -57747.09

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H113").Value = 584.5484
$ws_CUL.Range("I113").Value = 569.0833
$ws_CUL.Range("J113").Value = 594.3158
$ws_CUL.Range("K113").Value = 1707.2499
$ws_CUL.Range("L113").Value = 1782.9474
$ws_CUL.Range("M113").Value = 462.7501
$ws_CUL.Range("N113").Value = -6122.9474

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 25052.23
$ws_GSM.Range("I70").Value = 38500.516
$ws_GSM.Range("J70").Value = 5200
$ws_GSM.Range("K70").Value = 38500.516
$ws_GSM.Range("L70").Value = 5200
$ws_GSM.Range("M70").Value = -38230.516
$ws_GSM.Range("N70").Value = -5740

$ws_GSM.Range("H73").Value = 25052.23
$ws_GSM.Range("I73").Value = 38500.516
$ws_GSM.Range("J73").Value = 5200
$ws_GSM.Range("K73").Value = 38500.516
$ws_GSM.Range("L73").Value = 5200
$ws_GSM.Range("M73").Value = -37564.516
$ws_GSM.Range("N73").Value = -7072

$ws_GSM.Range("H97").Value = 5000
$ws_GSM.Range("I97").Value = 5000
$ws_GSM.Range("K97").Value = 5000
$ws_GSM.Range("M97").Value = -4504

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H109").Value = 28171.666
$ws_LTW.Range("J109").Value = 28171.666
$ws_LTW.Range("L109").Value = 28171.666
$ws_LTW.Range("N109").Value = -30945.666

$ws_LTW.Range("H132").Value = 2403.5652
$ws_LTW.Range("I132").Value = 1971.2727
$ws_LTW.Range("J132").Value = 2799.8333
$ws_LTW.Range("K132").Value = 5913.8181
$ws_LTW.Range("L132").Value = 8399.499899999999
$ws_LTW.Range("M132").Value = -3383.8181
$ws_LTW.Range("N132").Value = -13459.4999

$ws_LTW.Range("H133").Value = 52641.9
$ws_LTW.Range("J133").Value = 52641.9
$ws_LTW.Range("L133").Value = 52641.9
$ws_LTW.Range("N133").Value = -57701.9

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H109").Value = 29238.5
$ws_WVR.Range("J109").Value = 29238.5
$ws_WVR.Range("L109").Value = 29238.5
$ws_WVR.Range("N109").Value = -32012.5

$ws_WVR.Range("H132").Value = 2728.1794
$ws_WVR.Range("I132").Value = 2330.0334
$ws_WVR.Range("J132").Value = 4055.3333
$ws_WVR.Range("K132").Value = 6990.100199999999
$ws_WVR.Range("L132").Value = 12165.9999
$ws_WVR.Range("M132").Value = -4460.100199999999
$ws_WVR.Range("N132").Value = -17225.9999
